# Apply cryptocurrency price/volume updates from latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "70.882.94"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.532.47"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "613.33"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "174.16"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.30%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.527.89"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.74%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.611"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.23%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.04%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.46"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.24%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.589"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "46.73"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000277"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.23%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.104.24"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("E16").Value = "  +0.57%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "616.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.530.83"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.93%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "70.876.02"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.121"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.57%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.80"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.887"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.02"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.66%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "15.79"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "98.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -0.10%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "33.90"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.06%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  -0.76%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.18"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.25%  "
$ws.Range("E33").Value = "  -0.26%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.21%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "615.47"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.34%  "
$ws.Range("E36").Value = "  -0.33%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "10.87"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.55"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  +0.37%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "57.07"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +0.76%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.378.12"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0₃0742"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.38%  "
$ws.Range("E45").Value = "  -1.83%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "32.41"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.83%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  +0.10%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "134.16"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  -0.01%  "
